# The commit swaps the two triples of columns holding the "sib1" stats:
#   J:L  -> sib1_sex0_est / sib1_sex0_se / sib1_sex0_n
#   M:O  -> sib1_sex1_est / sib1_sex1_se / sib1_sex1_n
# After the edit, J:L holds what used to be in M:O (sib1_sex1_*) and M:O
# holds what used to be in J:L (sib1_sex0_*) - for every row, including the
# header row (so the column captions move together with their data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

$rangeLeft  = $ws.Range("J1:L$lastRow")
$rangeRight = $ws.Range("M1:O$lastRow")

$leftValues  = $rangeLeft.Value2
$rightValues = $rangeRight.Value2

$rangeLeft.Value2  = $rightValues
$rangeRight.Value2 = $leftValues
